$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "AG6304A74"
$ws.Range("A2").Value = "AG6304A47"
$ws.Range("A3").Value = "AG6304A79"
$ws.Range("A4").Value = "AG6304A100"
$ws.Range("A5").Value = "AG6304A117"

$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Select()
